# Facilitators guidelines - Conditional Probability.docx
# Translate the English UI/table labels (and one inline word) to Swahili.
#
# NOTE on ordering: "Video Introduction" is a substring of
# "General VMC Video Introduction", so the longer phrase must be replaced
# before the shorter one to avoid a partial / incorrect match.

$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute(
        $find,      # FindText
        $true,      # MatchCase
        $true,      # MatchWholeWord
        $false,     # MatchWildcards
        $false,     # MatchSoundsLike
        $false,     # MatchAllWordForms
        $true,      # Forward
        1,          # Wrap (wdFindContinue)
        $false,     # Format
        $replace,   # ReplaceWith
        2           # Replace (wdReplaceAll)
    ) | Out-Null
}

# Longer/more-specific phrase first to avoid being shadowed by the shorter one.
Replace-Text "General VMC Video Introduction" "Utangulizi Mkuu wa Video ya VMC"

Replace-Text "Video Title" "Kichwa cha Video"
Replace-Text "Topic" "Mada"
Replace-Text "Aim(s)" "Malengo"
Replace-Text "Length" "Urefu"
Replace-Text "Camp Location" "Mahali pa Kambi"
Replace-Text "Facilitators" "Wawezeshaji"
Replace-Text "N. of students" "N. ya wanafunzi"
Replace-Text "Date" "Tarehe"
Replace-Text "Resources" "Rasilimali"
Replace-Text "needed" "inahitajika"
Replace-Text "Preparations" "Maandalizi"
Replace-Text "Video time" "Muda wa video"
Replace-Text "What facilitator does" "Mwezeshaji anafanya nini"
Replace-Text "What learners do" "Wanachofanya wanafunzi"
Replace-Text "Video Introduction" "Utangulizi wa Video"

# Partial-sentence edit: only the "Why?" becomes "Mbona?"; rest is untouched.
Replace-Text `
    "Facilitate the discussion: the “dealer” never told anything about the winning card, and nonetheless the “guesser” guessed right more times after the hint than before. Why? How can this be viewed in terms of the first experiment?" `
    "Facilitate the discussion: the “dealer” never told anything about the winning card, and nonetheless the “guesser” guessed right more times after the hint than before. Mbona? How can this be viewed in terms of the first experiment?"

"done"
